# 5V External Buck Done (#7)
# Adds the "actual R_ON" input and the resulting t_ON calculation to the
# 5V buck sizing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New input: R_ON_ACT (row 18), formatted like the other scientific-notation
# inputs (K in B6, f_S in B7).
$ws.Range("A18").Value = "R_ON_ACT"
$ws.Range("B18").Value = 44200
$ws.Range("B18").NumberFormat = $ws.Range("B7").NumberFormat

# Register the defined name so formulas can refer to it by name.
$wb.Names.Add("R_ON_ACT", "=Sheet1!`$B`$18")

# New output: t_ON (row 19), computed from the newly added R_ON_ACT.
$ws.Range("A19").Value = "t_ON"
$ws.Range("B19").Formula = "=K*R_ON_ACT/V_IN"

# Move the active selection to the new last cell, matching the saved state.
$ws.Range("B19").Select() | Out-Null
